# Separate artist and album list in the wireframes
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "Artists" sheet before "Albums" and populate it with
#    the distinct list of artist names.
# ---------------------------------------------------------------------
$albumsSheet = $wb.Worksheets.Item("Albums")
$artists = $wb.Worksheets.Add($albumsSheet)
$artists.Name = "Artists"

$artists.Range("A1").Value = "Name"
$artists.Range("A1").Font.Bold = $true

$artists.Range("A2").Value = "Dire Straits"
$artists.Range("A3").Value = "Dragonette"
$artists.Range("A4").Value = "John Coltrane"
$artists.Range("A5").Value = "The Beatles"

$artists.Columns.Item(1).AutoFit() | Out-Null
$artists.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Trim the "Albums" sheet down to the two Beatles albums.
# ---------------------------------------------------------------------
$albums = $wb.Worksheets.Item("Albums")
$albums.Rows.Item(5).Delete()
$albums.Rows.Item(3).Delete()
$albums.Rows.Item(2).Delete()
$albums.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Replace the "Tracks" sheet contents with The Beatles (White Album)
#    track listing.
# ---------------------------------------------------------------------
$tracks = $wb.Worksheets.Item("Tracks")

$titles = @(
    "Back In The U.S.S.R.",
    "Dear Prudence",
    "Glass Onion",
    "Ob-La-Di, Ob-La-Da",
    "Wild Honey Pie",
    "The Continuing Story Of Bungalow Bill",
    "While My Guitar Gently Weeps",
    "Happiness Is A Warm Gun",
    "Martha My Dear",
    "I'm So Tired",
    "Blackbird",
    "Piggies",
    "Rocky Raccoon",
    "Don't Pass Me By",
    "Why Don't We Do It In The Road?",
    "I Will",
    "Julia"
)

$durations = @(
    0.11319444444444444,
    0.16388888888888889,
    0.095138888888888884,
    0.13055555555555556,
    0.036111111111111115,
    0.13472222222222222,
    0.19791666666666666,
    0.11319444444444444,
    0.10277777777777779,
    0.085416666666666655,
    0.095833333333333326,
    0.086111111111111124,
    0.14722222222222223,
    0.15972222222222224,
    0.07013888888888889,
    0.073611111111111113,
    0.12083333333333333
)

# Insert the extra rows needed (5 existing data rows -> 17 new data rows).
for ($i = 0; $i -lt 12; $i++) {
    $tracks.Rows.Item(7).Insert()
}

for ($i = 0; $i -lt 17; $i++) {
    $r = 2 + $i
    $tracks.Cells.Item($r, 1).Value = "The Beatles"
    $tracks.Cells.Item($r, 2).Value = "The Beatles"
    $tracks.Cells.Item($r, 3).Value = ($i + 1)
    $tracks.Cells.Item($r, 4).Value = $titles[$i]
    $tracks.Cells.Item($r, 5).Value = $durations[$i]
}

$tracks.Columns.Item(1).AutoFit() | Out-Null
$tracks.Columns.Item(4).AutoFit() | Out-Null
$tracks.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Tracks is the sheet left on-screen / active when the file was saved.
# ---------------------------------------------------------------------
$tracks.Activate()
